# Update weekly Fruta / Hortaliza price records (Pepino dulce) to reflect
# refreshed data rows, per upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44391
$ws.Range("I2").Value = 'Segunda'
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("P2").Value = 861

$ws.Range("D3").Value = 44433
$ws.Range("H3").Value = 'Cultivar IV Región'
$ws.Range("I3").Value = 'Segunda'
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("N3").Value = '$/bandeja 18 kilos'
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 972
$ws.Range("Q3").Value = 18

$ws.Range("D4").Value = 44433
$ws.Range("I4").Value = 'Tercera'
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 806

$ws.Range("D5").Value = 44405
$ws.Range("J5").Value = 140

$ws.Range("D6").Value = 44454
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 19000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 19500
$ws.Range("P6").Value = 1083

$ws.Range("D10").Value = 44412
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 150

$ws.Range("D11").Value = 44398
$ws.Range("I11").Value = 'Primera'

$ws.Range("D12").Value = 44398
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("P12").Value = 861

$ws.Range("D13").Value = 44363
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("P13").Value = 806

$ws.Range("D14").Value = 44221
$ws.Range("H14").Value = 'Cultivar XV región'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("N14").Value = '$/caja 10 kilos'
$ws.Range("O14").Value = 'Región de Arica y Parinacota'
$ws.Range("P14").Value = 550
$ws.Range("Q14").Value = 10

$ws.Range("D15").Value = 44435
$ws.Range("H15").Value = 'Cultivar IV Región'
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17500
$ws.Range("N15").Value = '$/bandeja 18 kilos'
$ws.Range("O15").Value = 'Provincia de Limarí'
$ws.Range("P15").Value = 972
$ws.Range("Q15").Value = 18

$ws.Range("D16").Value = 44435
$ws.Range("I16").Value = 'Tercera'
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("P16").Value = 806

$ws.Range("D17").Value = 44377
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17600
$ws.Range("P17").Value = 978

$ws.Range("D18").Value = 44211
$ws.Range("H18").Value = 'Cultivar XV región'
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = 4500
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 4750
$ws.Range("N18").Value = '$/caja 10 kilos'
$ws.Range("O18").Value = 'Región de Arica y Parinacota'
$ws.Range("P18").Value = 475
$ws.Range("Q18").Value = 10

